# Update the test data values on the "input" sheet.
# Rows 2-8 get new randomized string values (replacing the old ones,
# some reordered), and the last row (row 9) is cleared out so the
# sheet now only has 8 data values instead of 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(1, 1).Value = "Test Data"
$ws.Cells.Item(2, 1).Value = "etaeeigaehr"
$ws.Cells.Item(3, 1).Value = "rtrtrydmms"
$ws.Cells.Item(4, 1).Value = "eEeOUuoaRHRGSa"
$ws.Cells.Item(5, 1).Value = "iers-ta*_ta!d"
$ws.Cells.Item(6, 1).Value = "eouaeiuo"
$ws.Cells.Item(7, 1).Value = "e  gtr  w q  ii z "
$ws.Cells.Item(8, 1).Value = "1w1f1eg53qe4o"
$ws.Cells.Item(9, 1).Value = ""

# Reselect the last populated cell as the active selection.
$ws.Range("A9").Select()
